$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.656903666666667
$ws.Range("H2").Value = 13.970711
$ws.Range("I2").Value = 0.03472688986918189
$ws.Range("J2").Value = 0.03526504803992261
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 62.07563766666667
$ws.Range("N2").Value = 186.226913
$ws.Range("O2").Value = 0.4556169394345516
$ws.Range("P2").Value = 0.53808222397892
$ws.Range("Q2").Value = 289.0802646605715
$ws.Range("R2").Value = 2601.722381945143
$ws.Range("S2").Value = 0.01582215927827739
$ws.Range("T2").Value = 0.01897549547804501

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.656903666666667
$ws.Range("H3").Value = 13.970711
$ws.Range("I3").Value = 0.03472688986918189
$ws.Range("J3").Value = 0.03526504803992261
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.371752000000001
$ws.Range("N3").Value = 25.115256
$ws.Range("O3").Value = 0.06144619962548196
$ws.Range("P3").Value = 0.07256777544435759
$ws.Range("Q3").Value = 38.986442585224
$ws.Range("R3").Value = 350.877983267016
$ws.Range("S3").Value = 0.002133835407273878
$ws.Range("T3").Value = 0.002559106087195586

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.656903666666667
$ws.Range("H4").Value = 13.970711
$ws.Range("I4").Value = 0.03472688986918189
$ws.Range("J4").Value = 0.03526504803992261
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4371553333333333
$ws.Range("N4").Value = 1.311466
$ws.Range("O4").Value = 0.003208591687778628
$ws.Range("P4").Value = 0.003789337054374833
$ws.Range("Q4").Value = 2.035790274702889
$ws.Range("R4").Value = 18.322112472326
$ws.Range("S4").Value = 0.0001114244101766609
$ws.Range("T4").Value = 0.0001336311532619873

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.656903666666667
$ws.Range("H5").Value = 13.970711
$ws.Range("I5").Value = 0.03472688986918189
$ws.Range("J5").Value = 0.03526504803992261
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.718766
$ws.Range("N5").Value = 8.156298
$ws.Range("O5").Value = 0.01995494352567695
$ws.Range("P5").Value = 0.0235667277976885
$ws.Range("Q5").Value = 12.66103135420867
$ws.Range("R5").Value = 113.949282187878
$ws.Range("S5").Value = 0.0006929731261619276
$ws.Range("T5").Value = 0.0008310817879292644

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.656903666666667
$ws.Range("H6").Value = 13.970711
$ws.Range("I6").Value = 0.03472688986918189
$ws.Range("J6").Value = 0.03526504803992261
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 62.6419255
$ws.Range("N6").Value = 125.283851
$ws.Range("O6").Value = 0.4597733257265108
$ws.Range("P6").Value = 0.3619939357246589
$ws.Range("Q6").Value = 291.7174125480101
$ws.Range("R6").Value = 1750.304475288061
$ws.Range("S6").Value = 0.01596649764729203
$ws.Range("T6").Value = 0.01276573353349075

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 123.304606
$ws.Range("H7").Value = 369.913818
$ws.Range("I7").Value = 0.9194919584818978
$ws.Range("J7").Value = 0.9337412077596615
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 62.07563766666667
$ws.Range("N7").Value = 186.226913
$ws.Range("O7").Value = 0.4556169394345516
$ws.Range("P7").Value = 0.53808222397892
$ws.Range("Q7").Value = 7654.212044687092
$ws.Range("R7").Value = 68887.90840218384
$ws.Range("S7").Value = 0.4189361119582041
$ws.Range("T7").Value = 0.5024295456920815

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 123.304606
$ws.Range("H8").Value = 369.913818
$ws.Range("I8").Value = 0.9194919584818978
$ws.Range("J8").Value = 0.9337412077596615
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.371752000000001
$ws.Range("N8").Value = 25.115256
$ws.Range("O8").Value = 0.06144619962548196
$ws.Range("P8").Value = 0.07256777544435759
$ws.Range("Q8").Value = 1032.275581889712
$ws.Range("R8").Value = 9290.480237007409
$ws.Range("S8").Value = 0.05649928643490406
$ws.Range("T8").Value = 0.06775952228784636

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 123.304606
$ws.Range("H9").Value = 369.913818
$ws.Range("I9").Value = 0.9194919584818978
$ws.Range("J9").Value = 0.9337412077596615
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4371553333333333
$ws.Range("N9").Value = 1.311466
$ws.Range("O9").Value = 0.003208591687778628
$ws.Range("P9").Value = 0.003789337054374833
$ws.Range("Q9").Value = 53.90326613746533
$ws.Range("R9").Value = 485.129395237188
$ws.Range("S9").Value = 0.002950274254964309
$ws.Range("T9").Value = 0.003538260157760395

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 123.304606
$ws.Range("H10").Value = 369.913818
$ws.Range("I10").Value = 0.9194919584818978
$ws.Range("J10").Value = 0.9337412077596615
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.718766
$ws.Range("N10").Value = 8.156298
$ws.Range("O10").Value = 0.01995494352567695
$ws.Range("P10").Value = 0.0235667277976885
$ws.Range("Q10").Value = 335.236370436196
$ws.Range("R10").Value = 3017.127333925764
$ws.Range("S10").Value = 0.01834841010382036
$ws.Range("T10").Value = 0.02200522487675684

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 123.304606
$ws.Range("H11").Value = 369.913818
$ws.Range("I11").Value = 0.9194919584818978
$ws.Range("J11").Value = 0.9337412077596615
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 62.6419255
$ws.Range("N11").Value = 125.283851
$ws.Range("O11").Value = 0.4597733257265108
$ws.Range("P11").Value = 0.3619939357246589
$ws.Range("Q11").Value = 7724.037942858852
$ws.Range("R11").Value = 46344.22765715312
$ws.Range("S11").Value = 0.422757875730005
$ws.Range("T11").Value = 0.3380086547452163

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 6.139289
$ws.Range("H12").Value = 12.278578
$ws.Range("I12").Value = 0.04578115164892033
$ws.Range("J12").Value = 0.03099374420041592
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 62.07563766666667
$ws.Range("N12").Value = 186.226913
$ws.Range("O12").Value = 0.4556169394345516
$ws.Range("P12").Value = 0.53808222397892
$ws.Range("Q12").Value = 381.1002794949523
$ws.Range("R12").Value = 2286.601676969714
$ws.Range("S12").Value = 0.02085866819807016
$ws.Range("T12").Value = 0.01667718280879355

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 6.139289
$ws.Range("H13").Value = 12.278578
$ws.Range("I13").Value = 0.04578115164892033
$ws.Range("J13").Value = 0.03099374420041592
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.371752000000001
$ws.Range("N13").Value = 25.115256
$ws.Range("O13").Value = 0.06144619962548196
$ws.Range("P13").Value = 0.07256777544435759
$ws.Range("Q13").Value = 51.39660496432801
$ws.Range("R13").Value = 308.379629785968
$ws.Range("S13").Value = 0.002813077783304021
$ws.Range("T13").Value = 0.002249147069315643

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 6.139289
$ws.Range("H14").Value = 12.278578
$ws.Range("I14").Value = 0.04578115164892033
$ws.Range("J14").Value = 0.03099374420041592
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.4371553333333333
$ws.Range("N14").Value = 1.311466
$ws.Range("O14").Value = 0.003208591687778628
$ws.Range("P14").Value = 0.003789337054374833
$ws.Range("Q14").Value = 2.683822929224667
$ws.Range("R14").Value = 16.102937575348
$ws.Range("S14").Value = 0.0001468930226376586
$ws.Range("T14").Value = 0.0001174457433524511

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 6.139289
$ws.Range("H15").Value = 12.278578
$ws.Range("I15").Value = 0.04578115164892033
$ws.Range("J15").Value = 0.03099374420041592
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.718766
$ws.Range("N15").Value = 8.156298
$ws.Range("O15").Value = 0.01995494352567695
$ws.Range("P15").Value = 0.0235667277976885
$ws.Range("Q15").Value = 16.691290197374
$ws.Range("R15").Value = 100.147741184244
$ws.Range("S15").Value = 0.0009135602956946573
$ws.Range("T15").Value = 0.0007304211330023885

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 6.139289
$ws.Range("H16").Value = 12.278578
$ws.Range("I16").Value = 0.04578115164892033
$ws.Range("J16").Value = 0.03099374420041592
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 62.6419255
$ws.Range("N16").Value = 125.283851
$ws.Range("O16").Value = 0.4597733257265108
$ws.Range("P16").Value = 0.3619939357246589
$ws.Range("Q16").Value = 384.5768841609695
$ws.Range("R16").Value = 1750.304475288061
$ws.Range("S16").Value = 0.02104895234921383
$ws.Range("T16").Value = 0.01121954744595188
